# Rename the "About" diagram's cusp library reference.
# Slide 1, shape "Rounded Rectangle 25" (id=26) currently reads "cusplibrary";
# rename it to "cusparse".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $s.Shapes.Item("Rounded Rectangle 25")
$shape.TextFrame.TextRange.Text = "cusparse"
